$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J is column 10 (objective_value)

# --- Update objective_value results for newly-completed / corrected runs ---
$ws.Cells.Item(12, 10).Value = 2981.260000000134
$ws.Cells.Item(13, 10).Value = 3026.700000000139
$ws.Cells.Item(15, 10).Value = 3606.575000000219
$ws.Cells.Item(18, 10).Value = 4231.889999999888
$ws.Cells.Item(21, 10).Value = 4648.766666666827
$ws.Cells.Item(24, 10).Value = 6357.333333333332
$ws.Cells.Item(37, 10).Value = 6845.259999999817
$ws.Cells.Item(40, 10).Value = 8436.574999999741
$ws.Cells.Item(43, 10).Value = 10027.88999999892
$ws.Cells.Item(46, 10).Value = 11088.76666666706
$ws.Cells.Item(49, 10).Value = 13740.95833333339
$ws.Cells.Item(62, 10).Value = 14565.99333333231
$ws.Cells.Item(65, 10).Value = 18087.49166666673
$ws.Cells.Item(68, 10).Value = 21608.98999999871
$ws.Cells.Item(71, 10).Value = 23956.65555555516
$ws.Cells.Item(74, 10).Value = 29825.81944444447
$ws.Cells.Item(121, 10).Value = 2715.166666666654
$ws.Cells.Item(124, 10).Value = 3273.958333333345
$ws.Cells.Item(137, 10).Value = 3564.386666666736
$ws.Cells.Item(140, 10).Value = 4347.333333333303
$ws.Cells.Item(143, 10).Value = 5120.799999999786
$ws.Cells.Item(146, 10).Value = 5636.444444444604
$ws.Cells.Item(149, 10).Value = 6925.555555555547
$ws.Cells.Item(162, 10).Value = 7793.7933333333
$ws.Cells.Item(165, 10).Value = 9611.941666666642
$ws.Cells.Item(168, 10).Value = 11438.32999999957
$ws.Cells.Item(171, 10).Value = 12655.92222222231
$ws.Cells.Item(174, 10).Value = 15717.06944444447
$ws.Cells.Item(212, 10).Value = 2609.330000000022
$ws.Cells.Item(215, 10).Value = 3141.662500000004
$ws.Cells.Item(218, 10).Value = 3673.995000000001
$ws.Cells.Item(221, 10).Value = 4028.883333333339
$ws.Cells.Item(224, 10).Value = 4916.10416666667
$ws.Cells.Item(237, 10).Value = 5808.793492063502
$ws.Cells.Item(240, 10).Value = 7140.991865079373
$ws.Cells.Item(243, 10).Value = 8473.190238095194
$ws.Cells.Item(246, 10).Value = 9361.322486772508
$ws.Cells.Item(249, 10).Value = 11581.65310846563
$ws.Cells.Item(312, 10).Value = 7635.087142857112
$ws.Cells.Item(315, 10).Value = 9423.858928571421
$ws.Cells.Item(318, 10).Value = 11212.63071428553
$ws.Cells.Item(321, 10).Value = 12405.14523809527
$ws.Cells.Item(324, 10).Value = 15386.43154761906
$ws.Cells.Item(391, 10).Value = 91076.38826424249
$ws.Cells.Item(392, 10).Value = 91076.38826424249

# --- Clear objective_value for runs that no longer have a result ---
$ws.Cells.Item(362, 10).ClearContents()
$ws.Cells.Item(387, 10).ClearContents()
$ws.Cells.Item(399, 10).ClearContents()
$ws.Cells.Item(400, 10).ClearContents()
$ws.Cells.Item(401, 10).ClearContents()
